$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sample row was swapped out for a new student record (database di dir database)
$ws.Range("A2").Value = 9931231
$ws.Range("B2").Value = "Casgoni"
$ws.Range("C2").Value = "casgoni@gmail.com"

# Turn the e-mail cell into a real hyperlink (mailto:)
$ws.Hyperlinks.Add($ws.Range("C2"), "mailto:casgoni@gmail.com")
$ws.Range("C2").Style = "Hyperlink"

# Move the active selection onto the e-mail cell
$ws.Range("E2").Select()
